# Add a new column "CT" with data for 2024/12/15 (one day after the last
# existing date column "CS"). The header row stores the date as plain text
# (matching the style of the other date headers), and each data row stores a
# numeric value whose fill/style depends on its magnitude, mirroring the
# existing conditional-style convention used throughout the sheet:
#   style 1 -> value >= 140      (template cell A2, no fill)
#   style 2 -> value <  125      (template cell D2, yellow fill)
#   style 3 -> 125 <= value < 140 (template cell N2, light-blue fill)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell CT1: literal text "2024/12/15" (not a real date value) ---
$ws.Range("CT1").NumberFormat = "@"
$ws.Range("CT1").Value = "2024/12/15"
$ws.Range("A2").Copy()
$ws.Range("CT1").PasteSpecial(-4122)
$ws.Range("CT1").ColumnWidth = $ws.Range("CS1").ColumnWidth

# --- Data rows 2-53 ---
$newData = @(
    @{Row=2; Style=1; Value=145.4},
    @{Row=3; Style=1; Value=147.1},
    @{Row=4; Style=1; Value=145.8},
    @{Row=5; Style=1; Value=140.4},
    @{Row=6; Style=1; Value=163.6},
    @{Row=7; Style=2; Value=118.9},
    @{Row=8; Style=1; Value=161},
    @{Row=9; Style=1; Value=161.2},
    @{Row=10; Style=1; Value=177.6},
    @{Row=11; Style=3; Value=125.8},
    @{Row=12; Style=1; Value=149.8},
    @{Row=13; Style=1; Value=154.8},
    @{Row=14; Style=2; Value=117.4},
    @{Row=15; Style=2; Value=124.9},
    @{Row=16; Style=1; Value=190.7},
    @{Row=17; Style=1; Value=165.3},
    @{Row=18; Style=2; Value=124},
    @{Row=19; Style=1; Value=167.4},
    @{Row=20; Style=1; Value=185.9},
    @{Row=21; Style=1; Value=165.8},
    @{Row=22; Style=1; Value=157.7},
    @{Row=23; Style=3; Value=126.9},
    @{Row=24; Style=1; Value=191.5},
    @{Row=25; Style=2; Value=92.09999999999999},
    @{Row=26; Style=1; Value=159.2},
    @{Row=27; Style=2; Value=106.9},
    @{Row=28; Style=1; Value=146},
    @{Row=29; Style=1; Value=148},
    @{Row=30; Style=3; Value=134.2},
    @{Row=31; Style=1; Value=158},
    @{Row=32; Style=1; Value=148.6},
    @{Row=33; Style=1; Value=157.2},
    @{Row=34; Style=1; Value=165.1},
    @{Row=35; Style=3; Value=136.2},
    @{Row=36; Style=1; Value=161.3},
    @{Row=37; Style=1; Value=166.3},
    @{Row=38; Style=3; Value=127.6},
    @{Row=39; Style=1; Value=172.4},
    @{Row=40; Style=3; Value=131.5},
    @{Row=41; Style=3; Value=135.7},
    @{Row=42; Style=1; Value=174.4},
    @{Row=43; Style=1; Value=174.6},
    @{Row=44; Style=3; Value=131.5},
    @{Row=45; Style=1; Value=187},
    @{Row=46; Style=1; Value=156.3},
    @{Row=47; Style=1; Value=175.2},
    @{Row=48; Style=3; Value=128.4},
    @{Row=49; Style=1; Value=168.9},
    @{Row=50; Style=1; Value=156.3},
    @{Row=51; Style=1; Value=146.9},
    @{Row=52; Style=2; Value=118.6},
    @{Row=53; Style=2; Value=124.9}
)

$styleTemplates = @{1 = "A2"; 2 = "D2"; 3 = "N2"}

foreach ($item in $newData) {
    $cell = "CT" + $item.Row
    $template = $styleTemplates[[string]$item.Style]
    $ws.Range($template).Copy()
    $ws.Range($cell).PasteSpecial(-4122)
    $ws.Range($cell).Value = $item.Value
}
